# Scheduled market-data refresh: update currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) on a handful of rows across the ALC, ARM, BSM, CRP, CUL and GSM
# leve-profit sheets, mirroring the upstream price snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(53, 8).Value = 1646.8572
$ws.Cells.Item(53, 9).Value = 2504.6667
$ws.Cells.Item(53, 11).Value = 2504.6667
$ws.Cells.Item(53, 13).Value = -1867.6667
$ws.Cells.Item(98, 8).Value = 1445.1923
$ws.Cells.Item(98, 9).Value = 1278.2142
$ws.Cells.Item(98, 10).Value = 1640
$ws.Cells.Item(98, 11).Value = 1278.2142
$ws.Cells.Item(98, 12).Value = 1640
$ws.Cells.Item(98, 13).Value = 219.7858000000001
$ws.Cells.Item(98, 14).Value = -4636
$ws.Cells.Item(101, 8).Value = 1086.3077
$ws.Cells.Item(101, 9).Value = 739.7273
$ws.Cells.Item(101, 10).Value = 2992.5
$ws.Cells.Item(101, 11).Value = 2219.1819
$ws.Cells.Item(101, 12).Value = 8977.5
$ws.Cells.Item(101, 13).Value = -597.1819
$ws.Cells.Item(101, 14).Value = -12221.5
$ws.Cells.Item(122, 8).Value = 1445.1923
$ws.Cells.Item(122, 9).Value = 1278.2142
$ws.Cells.Item(122, 10).Value = 1640
$ws.Cells.Item(122, 11).Value = 3834.6426
$ws.Cells.Item(122, 12).Value = 4920
$ws.Cells.Item(122, 13).Value = -1384.6426
$ws.Cells.Item(122, 14).Value = -9820
$ws.Cells.Item(132, 8).Value = 2721.1887
$ws.Cells.Item(132, 9).Value = 2478.093
$ws.Cells.Item(132, 10).Value = 3766.5
$ws.Cells.Item(132, 11).Value = 7434.279
$ws.Cells.Item(132, 12).Value = 11299.5
$ws.Cells.Item(132, 13).Value = -4904.279
$ws.Cells.Item(132, 14).Value = -16359.5
$ws.Cells.Item(138, 8).Value = 2026.9166
$ws.Cells.Item(138, 9).Value = 1277.68
$ws.Cells.Item(138, 10).Value = 2841.3044
$ws.Cells.Item(138, 11).Value = 3833.04
$ws.Cells.Item(138, 12).Value = 8523.913199999999
$ws.Cells.Item(138, 13).Value = 1306.96
$ws.Cells.Item(138, 14).Value = -18803.9132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 15154798
$ws.Cells.Item(61, 9).Value = 20003214
$ws.Cells.Item(61, 10).Value = 3499.125
$ws.Cells.Item(61, 11).Value = 20003214
$ws.Cells.Item(61, 12).Value = 3499.125
$ws.Cells.Item(61, 13).Value = -20003002
$ws.Cells.Item(61, 14).Value = -3923.125
$ws.Cells.Item(74, 8).Value = 7354858
$ws.Cells.Item(74, 9).Value = 10001193
$ws.Cells.Item(74, 10).Value = 3927.2222
$ws.Cells.Item(74, 11).Value = 10001193
$ws.Cells.Item(74, 12).Value = 3927.2222
$ws.Cells.Item(74, 13).Value = -10000319
$ws.Cells.Item(74, 14).Value = -5675.2222
$ws.Cells.Item(77, 8).Value = 7354858
$ws.Cells.Item(77, 9).Value = 10001193
$ws.Cells.Item(77, 10).Value = 3927.2222
$ws.Cells.Item(77, 11).Value = 50005965
$ws.Cells.Item(77, 12).Value = 19636.111
$ws.Cells.Item(77, 13).Value = -50001597
$ws.Cells.Item(77, 14).Value = -28372.111
$ws.Cells.Item(97, 8).Value = 5474.5713
$ws.Cells.Item(97, 9).Value = 6593.2354
$ws.Cells.Item(97, 10).Value = 720.25
$ws.Cells.Item(97, 11).Value = 6593.2354
$ws.Cells.Item(97, 12).Value = 720.25
$ws.Cells.Item(97, 13).Value = -6097.2354
$ws.Cells.Item(97, 14).Value = -1712.25
$ws.Cells.Item(132, 8).Value = 7355154
$ws.Cells.Item(132, 9).Value = 12501934
$ws.Cells.Item(132, 10).Value = 2611.8572
$ws.Cells.Item(132, 11).Value = 37505802
$ws.Cells.Item(132, 12).Value = 7835.571599999999
$ws.Cells.Item(132, 13).Value = -37503272
$ws.Cells.Item(132, 14).Value = -12895.5716
$ws.Cells.Item(136, 8).Value = 15154798
$ws.Cells.Item(136, 9).Value = 20003214
$ws.Cells.Item(136, 10).Value = 3499.125
$ws.Cells.Item(136, 11).Value = 60009642
$ws.Cells.Item(136, 12).Value = 10497.375
$ws.Cells.Item(136, 13).Value = -60007092
$ws.Cells.Item(136, 14).Value = -15597.375
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1914.8948
$ws.Cells.Item(107, 9).Value = 1966.9286
$ws.Cells.Item(107, 10).Value = 1769.2
$ws.Cells.Item(107, 11).Value = 1966.9286
$ws.Cells.Item(107, 12).Value = 1769.2
$ws.Cells.Item(107, 13).Value = -46.92859999999996
$ws.Cells.Item(107, 14).Value = -5609.2
$ws.Cells.Item(132, 8).Value = 52078.25
$ws.Cells.Item(132, 10).Value = 52078.25
$ws.Cells.Item(132, 12).Value = 52078.25
$ws.Cells.Item(132, 14).Value = -62198.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 3375.9429
$ws.Cells.Item(58, 9).Value = 1116
$ws.Cells.Item(58, 10).Value = 6765.857
$ws.Cells.Item(58, 11).Value = 1116
$ws.Cells.Item(58, 12).Value = 6765.857
$ws.Cells.Item(58, 13).Value = -913
$ws.Cells.Item(58, 14).Value = -7171.857
$ws.Cells.Item(94, 8).Value = 9860.286
$ws.Cells.Item(94, 9).Value = 5005.5
$ws.Cells.Item(94, 10).Value = 16333.333
$ws.Cells.Item(94, 11).Value = 5005.5
$ws.Cells.Item(94, 12).Value = 16333.333
$ws.Cells.Item(94, 13).Value = -4554.5
$ws.Cells.Item(94, 14).Value = -17235.333
$ws.Cells.Item(99, 8).Value = 1443.8334
$ws.Cells.Item(99, 9).Value = 1351.5
$ws.Cells.Item(99, 11).Value = 1351.5
$ws.Cells.Item(99, 13).Value = 146.5
$ws.Cells.Item(126, 8).Value = 1443.8334
$ws.Cells.Item(126, 9).Value = 1351.5
$ws.Cells.Item(126, 11).Value = 4054.5
$ws.Cells.Item(126, 13).Value = -1584.5
$ws.Cells.Item(136, 8).Value = 3375.9429
$ws.Cells.Item(136, 9).Value = 1116
$ws.Cells.Item(136, 10).Value = 6765.857
$ws.Cells.Item(136, 11).Value = 3348
$ws.Cells.Item(136, 12).Value = 20297.571
$ws.Cells.Item(136, 13).Value = -798
$ws.Cells.Item(136, 14).Value = -25397.571
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 550.25
$ws.Cells.Item(17, 9).Value = 100.5
$ws.Cells.Item(17, 10).Value = 1000
$ws.Cells.Item(17, 11).Value = 301.5
$ws.Cells.Item(17, 12).Value = 3000
$ws.Cells.Item(17, 13).Value = -132.5
$ws.Cells.Item(17, 14).Value = -3338
$ws.Cells.Item(34, 8).Value = 1222.44
$ws.Cells.Item(34, 9).Value = 250
$ws.Cells.Item(34, 10).Value = 1307
$ws.Cells.Item(34, 11).Value = 750
$ws.Cells.Item(34, 12).Value = 3921
$ws.Cells.Item(34, 13).Value = -666
$ws.Cells.Item(34, 14).Value = -4089
$ws.Cells.Item(39, 8).Value = 522.0540999999999
$ws.Cells.Item(39, 10).Value = 522.0540999999999
$ws.Cells.Item(39, 12).Value = 1566.1623
$ws.Cells.Item(39, 14).Value = -2154.1623
$ws.Cells.Item(55, 8).Value = 550.2941
$ws.Cells.Item(55, 9).Value = 0
$ws.Cells.Item(55, 10).Value = 550.2941
$ws.Cells.Item(55, 11).Value = 0
$ws.Cells.Item(55, 12).Value = 1650.8823
$ws.Cells.Item(55, 13).ClearContents()
$ws.Cells.Item(55, 14).Value = -2004.8823
$ws.Cells.Item(131, 8).Value = 1005.5156
$ws.Cells.Item(131, 9).Value = 927.5
$ws.Cells.Item(131, 10).Value = 1016.6607
$ws.Cells.Item(131, 11).Value = 2782.5
$ws.Cells.Item(131, 12).Value = 3049.9821
$ws.Cells.Item(131, 13).Value = 2257.5
$ws.Cells.Item(131, 14).Value = -13129.9821
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 53731.156
$ws.Cells.Item(113, 9).Value = 59956
$ws.Cells.Item(113, 10).Value = 820
$ws.Cells.Item(113, 11).Value = 59956
$ws.Cells.Item(113, 12).Value = 820
$ws.Cells.Item(113, 13).Value = -57786
$ws.Cells.Item(113, 14).Value = -5160
$ws.Cells.Item(138, 8).Value = 57356.57
$ws.Cells.Item(138, 10).Value = 57356.57
$ws.Cells.Item(138, 12).Value = 57356.57
$ws.Cells.Item(138, 14).Value = -67636.57000000001
